$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (tab name) to match the new test-case id range.
$ws.Name = "SCD0011"

# Update the TC_ID cell (B2) to the new test case id.
$ws.Range("B2").Value = "SCD0011-020"

# Widen column B so the new, longer TC_ID value is fully visible (bestFit).
$ws.Columns("B").ColumnWidth = 12.43

# Select the whole header + data row and left-align / vertically center it
# (keeps each cell's existing font, number format and wrap settings).
$rng = $ws.Range("A1:Z2")
$rng.HorizontalAlignment = -4131
$rng.VerticalAlignment = -4108

# Move the active selection to B3, matching where the user left off editing.
$ws.Range("B3").Select()
